$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price record was inserted above the existing row 186 (now row 187),
# pushing every subsequent record down by one row.
$ws.Rows("186:186").Insert()

$ws.Range("A186").Value = 11
$ws.Range("B186").Value = "Vega Monumental Concepción"
$ws.Range("C186").Value = "Bíobío"
$ws.Range("D186").Value = 45205
$ws.Range("E186").Value = 8
$ws.Range("F186").Value = 100112032
$ws.Range("G186").Value = "Zapallo italiano"
$ws.Range("H186").Value = "Sin especificar"
$ws.Range("I186").Value = "Primera"
$ws.Range("J186").Value = 100
$ws.Range("K186").Value = 22000
$ws.Range("L186").Value = 23000
$ws.Range("M186").Value = 22500
$ws.Range("N186").Value = "$/caja 50 unidades"
$ws.Range("O186").Value = "Región de Arica y Parinacota"
$ws.Range("P186").Value = 450
$ws.Range("Q186").Value = 50
$ws.Range("R186").Value = "Hortaliza"
